$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure the Price (D) and Volume (E) columns are treated as plain text so
# values such as "27.205.35" or "16.50" are not reinterpreted as numbers/dates
# and so trailing zeros are preserved exactly as scraped.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.205.35"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.686.76"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.93"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.01"
$ws.Range("E8").Value = "  +13.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("E9").Value = "  +3.74%  "
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0891"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.925.53"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.692.82"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.18"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.553"
$ws.Range("E15").Value = "  +4.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.25"
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.205.81"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.28"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.17"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0745"
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.57"
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.65"
$ws.Range("E23").Value = "  +4.94%  "
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.22"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.31"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.50"
$ws.Range("E27").Value = "  +2.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.113"
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0502"
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.574.47"
$ws.Range("E33").Value = "  +6.31%  "
$ws.Range("E34").Value = "  +2.41%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.958"
$ws.Range("E36").Value = "  +5.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.603"
$ws.Range("E37").Value = "  +3.60%  "
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0174"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.06"
$ws.Range("E40").Value = "  +3.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.58"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.70"
$ws.Range("E43").Value = "  -2.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.26"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.833.06"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.786"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "91.27"
$ws.Range("E47").Value = "  +0.98%  "
$ws.Range("E48").Value = "  +5.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0109"
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("E50").Value = "  +3.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.20"
$ws.Range("E51").Value = "  +5.66%  "
